# Apply the "Updated cryptos list" data refresh (Mon May 22 19:49:24 UTC 2023).
# For every changed row we overwrite the B (Coin), C (Link), D (Price) and/or
# E (Volume 1h) cells with their new values. D-column values that look like a
# plain decimal number (e.g. "19.80", "0.8744") would otherwise be silently
# re-interpreted by Excel as a numeric value (losing the trailing zero /
# exact text formatting seen in the source data), so for those cells we force
# the cell to Text, write the value, then restore the default "Normal" style
# so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $value
    $rng.Style = 'Normal'
}

# Row 2
$ws.Range('D2').Value = '26.972.00'
$ws.Range('E2').Value = '  -0.48%  '

# Row 3
$ws.Range('D3').Value = '1.825.60'
$ws.Range('E3').Value = '  +0.12%  '

# Row 4
Set-TextValue 'D4' '1.004'
$ws.Range('E4').Value = '  -0.48%  '

# Row 5
Set-TextValue 'D5' '312.00'
$ws.Range('E5').Value = '  +0.14%  '

# Row 6
$ws.Range('E6').Value = '  -0.31%  '

# Row 7
Set-TextValue 'D7' '0.4624'
$ws.Range('E7').Value = '  -0.15%  '

# Row 8
$ws.Range('E8').Value = '  +1.81%  '

# Row 9
Set-TextValue 'D9' '0.07332'
$ws.Range('E9').Value = '  +0.62%  '

# Row 10
Set-TextValue 'D10' '0.8744'

# Row 11
Set-TextValue 'D11' '0.07928'
$ws.Range('E11').Value = '  +4.33%  '

# Row 12
Set-TextValue 'D12' '19.80'
$ws.Range('E12').Value = '  -1.62%  '

# Row 13
$ws.Range('D13').Value = '1.842.32'
$ws.Range('E13').Value = '  +0.54%  '

# Row 14
Set-TextValue 'D14' '5.335'
$ws.Range('E14').Value = '  -0.21%  '

# Row 15
Set-TextValue 'D15' '6.536'
$ws.Range('E15').Value = '  +0.96%  '

# Row 16
Set-TextValue 'D16' '91.18'
$ws.Range('E16').Value = '  -1.51%  '

# Row 17
$ws.Range('E17').Value = '  -0.12%  '

# Row 18
Set-TextValue 'D18' '0.000008864'
$ws.Range('E18').Value = '  +2.57%  '

# Row 19
Set-TextValue 'D19' '1.004'
$ws.Range('E19').Value = '  -0.41%  '

# Row 20
Set-TextValue 'D20' '14.77'
$ws.Range('E20').Value = '  +1.93%  '

# Row 21
$ws.Range('D21').Value = '27.007.62'
$ws.Range('E21').Value = '  -1.35%  '

# Row 22
Set-TextValue 'D22' '5.101'
$ws.Range('E22').Value = '  -1.88%  '

# Row 23
Set-TextValue 'D23' '10.54'
$ws.Range('E23').Value = '  -0.23%  '

# Row 24
$ws.Range('D24').Value = '2.042.91'
$ws.Range('E24').Value = '  -2.35%  '

# Row 25
Set-TextValue 'D25' '153.10'
$ws.Range('E25').Value = '  +0.91%  '

# Row 26
Set-TextValue 'D26' '1.844'
$ws.Range('E26').Value = '  -1.44%  '

# Row 28
Set-TextValue 'D28' '2.037'
$ws.Range('E28').Value = '  -3.04%  '

# Row 29
Set-TextValue 'D29' '5.135'
$ws.Range('E29').Value = '  +1.09%  '

# Row 30
Set-TextValue 'D30' '115.46'
$ws.Range('E30').Value = '  -0.72%  '

# Row 31
Set-TextValue 'D31' '0.08900'
$ws.Range('E31').Value = '  -0.26%  '

# Row 32
Set-TextValue 'D32' '2.965'
$ws.Range('E32').Value = '  +0.18%  '

# Row 33
Set-TextValue 'D33' '0.7285'
$ws.Range('E33').Value = '  -0.96%  '

# Row 34
$ws.Range('E34').Value = '  -0.33%  '

# Row 35
$ws.Range('E35').Value = '  -0.93%  '

# Row 36
Set-TextValue 'D36' '2.490'
$ws.Range('E36').Value = '  -2.32%  '

# Row 37
Set-TextValue 'D37' '0.01952'
$ws.Range('E37').Value = '  +1.78%  '

# Row 38
$ws.Range('E38').Value = '  -0.36%  '

# Row 39
Set-TextValue 'D39' '0.05222'
$ws.Range('E39').Value = '  -0.75%  '

# Row 40
Set-TextValue 'D40' '2.938'
$ws.Range('E40').Value = '  +0.09%  '

# Row 41
Set-TextValue 'D41' '7.098'
$ws.Range('E41').Value = '  -1.01%  '

# Row 42
Set-TextValue 'D42' '0.5155'
$ws.Range('E42').Value = '  -1.13%  '

# Row 43
Set-TextValue 'D43' '0.1621'
$ws.Range('E43').Value = '  -0.75%  '

# Row 44
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D44' '8.174'
$ws.Range('E44').Value = '  -1.20%  '

# Row 45
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D45' '0.4844'
$ws.Range('E45').Value = '  -0.97%  '

# Row 46
Set-TextValue 'D46' '1.005'

# Row 47
Set-TextValue 'D47' '10.21'
$ws.Range('E47').Value = '  +0.64%  '

# Row 48
Set-TextValue 'D48' '102.76'
$ws.Range('E48').Value = '  -1.19%  '

# Row 50
Set-TextValue 'D50' '0.06192'
$ws.Range('E50').Value = '  -0.98%  '

# Row 51
Set-TextValue 'D51' '64.83'
$ws.Range('E51').Value = '  +0.45%  '
